$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("R1")
$ws.Range("G2").Value = "3928:47:28"
$ws.Range("G3").Value = "68:20:06"

$ws = $wb.Worksheets.Item("R2")
$ws.Range("G2").Value = "12110:11:09"
$ws.Range("G3").Value = "3239:54:38"
$ws.Range("G4").Value = "478:06:12"

$ws = $wb.Worksheets.Item("R4")
$ws.Range("G2").Value = "2956:00:58"
$ws.Range("G3").Value = "183:13:13"
$ws.Range("G4").Value = "71:25:38"
$ws.Range("G5").Value = "69:03:11"

$ws = $wb.Worksheets.Item("R5")
$ws.Range("G2").Value = "429:59:57"

$ws = $wb.Worksheets.Item("R6")
$ws.Range("G2").Value = "70:32:15"
